$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (PLAYER_Data headers): replace DirX/DirY/DirZ with Yaw/Pitch, keep currentHealth,
# replace maxHealth with currentWeaponID
$ws.Range("D8").Value = "Yaw"
$ws.Range("E8").Value = "Pitch"
$ws.Range("F8").Value = "currentHealth"

# New Weapon section
$ws.Range("A21").Value = "Weapon"
$ws.Range("C21").Value = "WEAPON_Data"

$ws.Range("G8").Value = "currentWeaponID"
$ws.Range("H8").ClearContents()

$ws.Range("A22").Value = "weaponID"
$ws.Range("B22").Value = "name"
$ws.Range("C22").Value = "damage"
$ws.Range("D22").Value = "interval"
$ws.Range("E22").Value = "range"

# Selection moves to F27
$ws.Range("F27").Select()
